# New Invoice GUI now correctly displays. Onto connecting the buttons and
# settings for it.
#
# Appends two new log-book entries (rows 65-66) to Sheet1 and fills in the
# previously-blank "end time" / notes cells of row 64, matching the next
# activity-summary entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 64: this entry now has a logged end-time and an extra note ---
$ws.Range("C64").Value = "4:34PM"
$ws.Range("C64").NumberFormat = "h:mm"
$ws.Range("H64").Value = "Debugged as well"

# --- Row 65: new log entry ---
$ws.Range("B65").Value = 0.19027777777777777
$ws.Range("B65").NumberFormat = "h:mm"
$ws.Range("F65").Value = "Code/Design"
$ws.Range("G65").Value = "Build Invoice GUI"
$ws.Range("C65").Value = "5:32PM"
$ws.Range("C65").NumberFormat = "h:mm"
$ws.Range("E65").Value = 58

# --- Row 66: new log entry ---
$ws.Range("B66").Value = "5:32PM"
$ws.Range("F66").Value = "Debug"
$ws.Range("G66").Value = "Debug GUI "

# --- Scroll/selection so the newly-added rows are in view ---
$ws.Range("G66").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 2
